$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 4 ("2.") Name cell: "dark ground" -> "Light  ground"
$cell1 = $t.Cell(4, 5)
$cell1.Range.Text = "Light  ground"

# Row 5 ("3.") Name cell: "light ground" -> "Dark  ground"
$cell2 = $t.Cell(5, 5)
$cell2.Range.Text = "Dark  ground"
